$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "I10" as the new C4 value, shifting the existing C4:C13 values down
# by one row, and re-append the final "I10" value at the new row 14.
# Also populate A14/B14 to match the existing pattern (P / ACS1001).

$values = @("I10","I31.0","T81.2","S25.0","Y60.0","Y92.22","J96.09","F17.1","I48.9","E87.6","I10")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 1).Value = "P"
    $ws.Cells.Item($row, 2).Value = "ACS1001"
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

$ws.Range("C4").Select()
